{"js": "// Update \"M\u00e9thode de travail\" scope statement: add a paragraph describing\n// the git flow versioning approach, with surrounding spacer paragraphs that\n// match the document's existing visual rhythm, and fix a small typo later\n// in the document (\"Un Outil\" -> \"Un outil\").\n\nconst GIT_FLOW_TEXT =\n  \"Pour g\u00e9rer le versionning nous utilisons la m\u00e9thodologie git flow. \" +\n  \"Nous partons d\\u2019une branche main qui contient le code finalis\u00e9. \" +\n  \"A partir de celle-ci cr\u00e9ons une branche par t\u00e2che de travail. \" +\n  \"Lorsque nous avons termin\u00e9 la fonctionnalit\u00e9 du code nous faisons un merge dans la branche main. \" +\n  \"Ceci nous permet de travailler de fa\u00e7on efficace et structur\u00e9e.\";\n\n// 1) Locate the empty paragraph that immediately precedes the paragraph that\n//    starts with \"Nous utiliserons la convention de commit\". In the source\n//    document this paragraph is blank (it only acts as a spacer).\nlet paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet commitConventionIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(\"Nous utiliserons la convention de commit\") !== -1) {\n    commitConventionIndex = i;\n    break;\n  }\n}\n\nif (commitConventionIndex === -1) {\n  throw new Error(\"Could not find the 'Nous utiliserons la convention de commit' paragraph.\");\n}\n\nconst spacerParagraph = paragraphs.items[commitConventionIndex - 1];\nspacerParagraph.load(\"text\");\nawait context.sync();\n\n// 2) Fill that spacer paragraph with the new git flow text, then insert a\n//    fresh blank paragraph right after it (before the commit-convention\n//    paragraph), so the new paragraph is followed by one empty line exactly\n//    like the rest of this section.\nspacerParagraph.insertText(GIT_FLOW_TEXT, \"Replace\");\nspacerParagraph.insertParagraph(\"\", \"After\");\nawait context.sync();\n\n// 3) Re-fetch paragraphs (indices shifted because of the insertion above)\n//    and add a blank paragraph right after the commit-convention paragraph\n//    (the one ending in \"...r\u00e9alis\u00e9s durant ce projet.\").\nparagraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet conventionParagraph = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const t = paragraphs.items[i].text;\n  if (t.indexOf(\"qui permet d\\u2019avoir une vue globale\") !== -1) {\n    conventionParagraph = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!conventionParagraph) {\n  throw new Error(\"Could not find the paragraph ending the commit convention sentence.\");\n}\n\nconventionParagraph.insertParagraph(\"\", \"After\");\nawait context.sync();\n\n// 4) Fix the small typo later in the document: \"Un Outil\" -> \"Un outil\"\n//    (in the DB diagram sentence).\nconst typoMatches = context.document.body.search(\"Un Outil\", { matchCase: true, matchWholeWord: false });\ntypoMatches.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < typoMatches.items.length; i++) {\n  typoMatches.items[i].insertText(\"Un outil\", \"Replace\");\n}\nawait context.sync();\n", "ps1": "# Update \"M\u00e9thode de travail\" scope statement: add a paragraph describing\n# the git flow versioning approach, with surrounding spacer paragraphs that\n# match the document's existing visual rhythm, and fix a small typo later\n# in the document (\"Un Outil\" -> \"Un outil\").\n\n$d = $word.ActiveDocument\n\n$gitFlowText = \"Pour g\u00e9rer le versionning nous utilisons la m\u00e9thodologie git flow. \" + `\n    \"Nous partons d\u2019une branche main qui contient le code finalis\u00e9. \" + `\n    \"A partir de celle-ci cr\u00e9ons une branche par t\u00e2che de travail. \" + `\n    \"Lorsque nous avons termin\u00e9 la fonctionnalit\u00e9 du code nous faisons un merge dans la branche main. \" + `\n    \"Ceci nous permet de travailler de fa\u00e7on efficace et structur\u00e9e.\"\n\n# 1) Locate the paragraph that starts the commit-convention sentence, then\n#    step back one paragraph to find the blank spacer paragraph right\n#    before it.\n$commitConventionIndex = -1\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    if ($d.Paragraphs.Item($i).Range.Text -like \"*Nous utiliserons la convention de commit*\") {\n        $commitConventionIndex = $i\n        break\n    }\n}\n\nif ($commitConventionIndex -eq -1) {\n    throw \"Could not find the 'Nous utiliserons la convention de commit' paragraph.\"\n}\n\n$spacerParagraph = $d.Paragraphs.Item($commitConventionIndex - 1)\n\n# 2) Fill that spacer paragraph with the new git flow text, then insert a\n#    fresh blank paragraph right after it (before the commit-convention\n#    paragraph), so the new paragraph is followed by one empty line exactly\n#    like the rest of this section.\n$spacerRange = $spacerParagraph.Range\n$spacerRange.Text = $gitFlowText\n$spacerRange.InsertParagraphAfter()\n\n# 3) Re-locate the paragraph that ends the commit-convention sentence and\n#    add a blank paragraph right after it.\n$conventionIndex = -1\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    if ($d.Paragraphs.Item($i).Range.Text -like \"*qui permet d\u2019avoir une vue globale*\") {\n        $conventionIndex = $i\n        break\n    }\n}\n\nif ($conventionIndex -eq -1) {\n    throw \"Could not find the paragraph ending the commit convention sentence.\"\n}\n\n$conventionParagraph = $d.Paragraphs.Item($conventionIndex)\n$conventionParagraph.Range.InsertParagraphAfter()\n\n# 4) Fix the small typo later in the document: \"Un Outil\" -> \"Un outil\"\n#    (in the DB diagram sentence). Use MatchCase so the already-correct\n#    lowercase occurrence elsewhere in the document is left untouched.\n$findRange = $d.Content\n$findRange.Find.ClearFormatting()\n$findRange.Find.MatchCase = $true\n$findRange.Find.MatchWholeWord = $false\n$found = $findRange.Find.Execute(\"Un Outil\")\nif ($found) {\n    $findRange.Text = \"Un outil\"\n}\n"}
